$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: B6, D6, G6, L6, M6
$ws.Range("B6").Value = "-"
$ws.Range("D6").Value = "169-Press Approval Task "
$ws.Range("G4").Copy()
$ws.Range("G6").PasteSpecial(-4163)
$ws.Range("L6").Value = "Press Approval Task"
$ws.Range("M6").Value = "169-Press Approval Task "

# Row 7: B7, D7, G7, L7, M7
$ws.Range("B7").Value = "Digital Print F 4x0"
$ws.Range("D7").Value = "252-HP 10000 Press"
$ws.Range("G2").Copy()
$ws.Range("G7").PasteSpecial(-4163)
$ws.Range("L7").Value = "252-HP 10000 Press"
$ws.Range("M7").Value = "252-HP 10000 Press"

# Row 8: B8, D8, G8, L8, M8
$ws.Range("B8").Value = "Cut"
$ws.Range("D8").Value = "406-45`" Polar 115ED Cutter"
$ws.Range("G3").Copy()
$ws.Range("G8").PasteSpecial(-4163)
$ws.Range("L8").Value = "406-45`" Polar 115ED Cutter"
$ws.Range("M8").Value = "406-45`" Polar 115ED Cutter`n404-45`" Polar 115EMC Cutter`n405-54`" Polar 137EMC Cutter`n402-45`" Polar 115EMC Cutter`n403-54`" Polar 137ED Cutter"
